$wb = $excel.ActiveWorkbook

# --- Sheet renames -----------------------------------------------------
$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Name = "data statistics summary"

$wsResults = $wb.Worksheets.Item("Results Table")
$wsResults.Name = "Misinf detection results table"

# --- TF-IDF sheet: add per-class accuracy column (contingency matrix) --
$ws = $wb.Worksheets.Item("TF-IDF")

# Narrow the old "total" column and size the new accuracy column
$ws.Columns.Item(7).ColumnWidth = 16.2
$ws.Columns.Item(8).ColumnWidth = 15.1

# Relabel the totals header to reflect it now sits next to per-class accuracy
$ws.Range("G1").Value = "Test items per class "

# New header for the accuracy column
$ws.Range("H1").Value = "Accuracy per class"
$ws.Range("H1").Interior.ThemeColor = 10
$ws.Range("H1").Interior.TintAndShade = 0.6

# Per-class accuracy: correctly predicted count (diagonal) / items in that class
$ws.Range("H2").Formula = "=B2/G2"
$ws.Range("H3").Formula = "=C3/G3"
$ws.Range("H4").Formula = "=D4/G4"
$ws.Range("H5").Formula = "=E5/G5"
$ws.Range("H6").Formula = "=F6/G6"

"H2","H3","H4","H5","H6" | ForEach-Object {
    $cell = $ws.Range($_)
    $cell.Interior.ThemeColor = 10
    $cell.Interior.TintAndShade = 0.6
}

# Overall/average accuracy across classes
$ws.Range("H7").Formula = "=AVERAGE(H2:H6)"
$ws.Range("H7").Interior.ThemeColor = 10
$ws.Range("H7").Interior.TintAndShade = 0.4

# Restore selection to match the post-edit cursor position
$ws.Range("D17").Select()
